$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D10 value to 100.0 (numeric), matching C10's pattern
$ws.Range("D10").Value = 100.0

# Collapse all defined columns (mark "collapsed" outline attribute = true)
$ws.Columns.Item(1).OutlineLevel = 1
for ($col = 1; $col -le 16384; $col++) {
    $ws.Columns.Item($col).OutlineLevel = 1
}
